$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the data of row 2 and row 4 (these two records were
# reordered), including moving the "Publik kommentar" value (column AC)
# from row 2 to row 4. Only the columns that actually differ between the
# two rows are touched, to avoid disturbing other (identical) cells.

$cols = @("A","B","D","E","F","G","H","Q","R","Z","AB")

foreach ($col in $cols) {
    $addr2 = "$col" + "2"
    $addr4 = "$col" + "4"
    $v2 = $ws.Range($addr2).Value2
    $v4 = $ws.Range($addr4).Value2
    $ws.Range($addr2).Value2 = $v4
    $ws.Range($addr4).Value2 = $v2
}

# Move the "Publik kommentar" (AC) text from row 2 to row 4.
$ac2 = $ws.Range("AC2").Value2
$ws.Range("AC4").Value2 = $ac2
$ws.Range("AC2").Value2 = $null
